$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 30
$ws.Range("D3").Value = 1

$ws.Range("B4").Value = 40
$ws.Range("C4").Value = 30
$ws.Range("D4").Value = 1

$ws.Range("B5").Value = 20
$ws.Range("C5").Value = 30
$ws.Range("D5").Value = 1

$lo = $ws.ListObjects.Item(1)
$col = $lo.ListColumns.Item(5)
$col.DataBodyRange.Formula = "=(B2/C2)*D2"

$ws.Range("C6").Select()
